# Generate Report for Handback
# Update the "Generate Date"/"Datetime" timestamp strings to reflect the
# latest handback report generation run. These cells are stored as plain
# text (not Excel date serials), so we explicitly write them as strings.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet - "Latest HO Xliff Generate Date" (column G)
$wsOverview.Range("G2").Value = "2016-08-19 19:13:17"

# zh-cn sheet - "Correspond Handoff Datetime" (H) and
#               "Correspond Handback DateTime" (K)
$wsZhCn.Range("H2").Value = "2016-08-19 19:13:13"
$wsZhCn.Range("K2").Value = "2016-08-19 19:13:29"

# de-de sheet - "Correspond Handoff Datetime" (H) and
#               "Correspond Handback DateTime" (K)
$wsDeDe.Range("H2").Value = "2016-08-19 19:13:17"
$wsDeDe.Range("K2").Value = "2016-08-19 19:13:35"
